# "Generate Report for handoff"
#
# The previous handoff record (GUID 8b56ef7e-...) is replaced by a fresh
# one (GUID 44f9dfa4-...) with an updated content hash and new
# handoff timestamps. The row that tracked the now-superseded
# "8ecb5a77-...md / Handoff transform failed / Ignored" record is removed
# entirely (it was the last data row on every sheet), and the
# ".localization-config" row moves up from row 4 to row 3 (i.e. the
# now-unused row 4 is dropped).

$wb = $excel.ActiveWorkbook

$oldGuid = "8b56ef7e-f9a1-4f63-ac56-9846e2a326bd"
$newGuid = "44f9dfa4-9b64-4e46-ae02-f2a609207392"
$oldHash = "260c43e2277132eb1bed5af25732cd9a607a2e47"
$newHash = "2c45da9cd283941d9a9ceaac624915d452b3895b"

$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/d28b1fd7376c1f3ee620e565b7654176ed53688f/e2e/"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/d28b1fd7376c1f3ee620e565b7654176ed53688f/.localization-config"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e5d5d529c42747faf11a1c769b1008c287815fb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/" + $newGuid + "." + $newHash + ".zh-cn.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/df913b8e4210cf69a52c29a1fd32d4790b13a8db/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/" + $newGuid + "." + $newHash + ".de-de.xlf"

$newMdName = $newGuid + ".md"
$newZhXlfName = $newGuid + "." + $newHash + ".zh-cn.xlf"
$newDeXlfName = $newGuid + "." + $newHash + ".de-de.xlf"

function Fix-Hyperlink($ws, $cellAddr, $address, $display) {
    $hls = $ws.Hyperlinks
    foreach ($h in $hls) {
        if ($h.Range.Address() -eq $cellAddr) {
            $h.Address = $address
            $h.TextToDisplay = $display
            return
        }
    }
}

function Delete-Hyperlink($ws, $cellAddr) {
    $hls = $ws.Hyperlinks
    foreach ($h in $hls) {
        if ($h.Range.Address() -eq $cellAddr) {
            $h.Delete()
            return
        }
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview": columns File Name / zh-cn / de-de
# ---------------------------------------------------------------------
$wsO = $wb.Worksheets.Item("Overview")

# Row 3 currently tracks the "8ecb5a77...md / Handoff transform failed"
# record. Repurpose its hyperlink in place to become the
# ".localization-config" row (what is currently row 4) before the row
# shift, so the single-cell "ref" formatting on the hyperlink is preserved.
Fix-Hyperlink $wsO '$A$3' $cfgUrl ".localization-config"
Delete-Hyperlink $wsO '$A$4'
Fix-Hyperlink $wsO '$A$2' ($mdBase + $newMdName) $newMdName

$wsO.Range("A2").Value = $newMdName
$wsO.Range("B2").Value = "Ready for handoff"
$wsO.Range("C2").Value = "Ready for handoff"

# Drop the now-superseded row (was row 3); row 4 (.localization-config)
# shifts up to become the new row 3.
$wsO.Rows.Item(3).Delete()

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

Fix-Hyperlink $wsZh '$A$3' $cfgUrl ".localization-config"
Delete-Hyperlink $wsZh '$A$4'
Fix-Hyperlink $wsZh '$A$2' ($mdBase + $newMdName) $newMdName
Fix-Hyperlink $wsZh '$C$2' $zhXlfUrl $newZhXlfName

$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("B2").Value = "Ready for handoff"
$wsZh.Range("C2").Value = $newZhXlfName
$wsZh.Range("D2").Value = "2016-01-13 12:53:34"
$wsZh.Range("G2").Value = "0001-01-01 00:00:00"
$wsZh.Range("H2").Value = "Include"

$wsZh.Rows.Item(3).Delete()

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

Fix-Hyperlink $wsDe '$A$3' $cfgUrl ".localization-config"
Delete-Hyperlink $wsDe '$A$4'
Fix-Hyperlink $wsDe '$A$2' ($mdBase + $newMdName) $newMdName
Fix-Hyperlink $wsDe '$C$2' $deXlfUrl $newDeXlfName

$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("B2").Value = "Ready for handoff"
$wsDe.Range("C2").Value = $newDeXlfName
$wsDe.Range("D2").Value = "2016-01-13 12:53:43"
$wsDe.Range("G2").Value = "0001-01-01 00:00:00"
$wsDe.Range("H2").Value = "Include"

$wsDe.Rows.Item(3).Delete()
